# Update carbon stats on row 2 of the tot_c_hist sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C2"  = 21
    "D2"  = 203
    "E2"  = 21527
    "F2"  = 751926
    "G2"  = 2390803
    "H2"  = 1968011
    "I2"  = 2358895
    "J2"  = 3719569
    "K2"  = 6156445
    "L2"  = 9174635
    "M2"  = 10463968
    "N2"  = 11342819
    "O2"  = 10024836
    "P2"  = 8602668
    "Q2"  = 7463163
    "R2"  = 8142096
    "S2"  = 7829014
    "T2"  = 7723568
    "U2"  = 8568763
    "V2"  = 9302383
    "W2"  = 8766615
    "X2"  = 7936589
    "Y2"  = 7730442
    "Z2"  = 5776685
    "AA2" = 4401692
    "AB2" = 3101886
    "AC2" = 1841453
    "AD2" = 1005191
    "AE2" = 492742
    "AF2" = 222080
    "AG2" = 104512
    "AH2" = 66768
    "AI2" = 40040
    "AJ2" = 15073
    "AK2" = 4456
    "AL2" = 1608
    "AM2" = 485
    "AN2" = 147
    "AO2" = 42
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
